$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 312; this shifts the existing rows 312:389
# down to 313:390 (old row 389 becomes new row 390).
$ws.Rows("312:312").Insert()

# Populate the newly inserted row 312 with its data. The shared columns
# (A,B,C,E,F,G,H,N,O,Q,R) match the same values as the surrounding rows.
$ws.Range("A312").Value = 7
$ws.Range("B312").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C312").Value = "Ñuble"
$ws.Range("D312").Value = 45211
$ws.Range("E312").Value = 16
$ws.Range("F312").Value = 100112045
$ws.Range("G312").Value = "Zapallo"
$ws.Range("H312").Value = "Camote"
$ws.Range("I312").Value = "1a nueva(o)"
$ws.Range("J312").Value = 200
$ws.Range("K312").Value = 800
$ws.Range("L312").Value = 800
$ws.Range("M312").Value = 800
$ws.Range("N312").Value = '$/kilo (volumen en unidades)'
$ws.Range("O312").Value = "Región de O'Higgins"
$ws.Range("P312").Value = 800
$ws.Range("Q312").Value = 1
$ws.Range("R312").Value = "Hortaliza"
